# Rename the header cells from the "_old"/"_new" suffix convention to the
# "_FV2310"/"_FV2404" format-version suffix convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerMap = @{
    "Segmentname_old"          = "Segmentname_FV2310";
    "Segmentgruppe_old"        = "Segmentgruppe_FV2310";
    "Segment_old"              = "Segment_FV2310";
    "Datenelement_old"         = "Datenelement_FV2310";
    "Segment ID_old"           = "Segment ID_FV2310";
    "Code_old"                 = "Code_FV2310";
    "Qualifier_old"            = "Qualifier_FV2310";
    "Beschreibung_old"         = "Beschreibung_FV2310";
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2310";
    "Bedingung_old"            = "Bedingung_FV2310";
    "Segmentname_new"          = "Segmentname_FV2404";
    "Segmentgruppe_new"        = "Segmentgruppe_FV2404";
    "Segment_new"              = "Segment_FV2404";
    "Datenelement_new"         = "Datenelement_FV2404";
    "Segment ID_new"           = "Segment ID_FV2404";
    "Code_new"                 = "Code_FV2404";
    "Qualifier_new"            = "Qualifier_FV2404";
    "Beschreibung_new"         = "Beschreibung_FV2404";
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2404";
    "Bedingung_new"            = "Bedingung_FV2404";
}

$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $current = $cell.Value()
    if ($headerMap.ContainsKey($current)) {
        $cell.Value = $headerMap[$current]
    }
}

# Turn the header row + data into a proper Excel Table (ListObject) so the
# column headers above are wired up as structured-table column names and an
# autofilter is shown.
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row (split below row 1, freeze panes).
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)

Write-Host "Header renaming, table creation and freeze panes applied."
